$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.45"
$ws.Range("G2").Value = "'21"
$ws.Range("G3").Value = "'21"
$ws.Range("D4").Value = "'5.418"
$ws.Range("G4").Value = "'21"
$ws.Range("D5").Value = "'0.05893"
$ws.Range("G5").Value = "'21"
$ws.Range("D6").Value = "'3.432"
$ws.Range("G6").Value = "'21"
$ws.Range("D7").Value = "'6.521"
$ws.Range("G7").Value = "'21"
$ws.Range("D8").Value = "'0.8078"
$ws.Range("G8").Value = "'21"
$ws.Range("D9").Value = "'0.9306"
$ws.Range("G9").Value = "'21"
$ws.Range("G10").Value = "'21"
$ws.Range("D11").Value = "'0.07393"
$ws.Range("G11").Value = "'21"
$ws.Range("D12").Value = "'0.03318"
$ws.Range("G12").Value = "'21"
$ws.Range("D13").Value = "'0.03072"
$ws.Range("G13").Value = "'21"
$ws.Range("D14").Value = "'0.09348"
$ws.Range("G14").Value = "'21"
$ws.Range("D15").Value = "'3.853"
$ws.Range("G15").Value = "'21"
$ws.Range("D16").Value = "'0.001576"
$ws.Range("G16").Value = "'21"
$ws.Range("D17").Value = "'0.04671"
$ws.Range("G17").Value = "'21"
$ws.Range("D18").Value = "'0.0005919"
$ws.Range("G18").Value = "'21"
$ws.Range("D19").Value = "'0.006000"
$ws.Range("G19").Value = "'21"
$ws.Range("E20").Value = "'19BitKanKAN"
$ws.Range("G20").Value = "'21"
$ws.Range("D21").Value = "'0.004900"
$ws.Range("G21").Value = "'21"
$ws.Range("D22").Value = "'0.00006799"
$ws.Range("G22").Value = "'21"
$ws.Range("D23").Value = "'3.563"
$ws.Range("G23").Value = "'21"
$ws.Range("D24").Value = "'2.144"
$ws.Range("G24").Value = "'21"
$ws.Range("G25").Value = "'21"
$ws.Range("D26").Value = "'0.1331"
$ws.Range("G26").Value = "'21"
$ws.Range("D27").Value = "'0.0002295"
$ws.Range("G27").Value = "'21"
$ws.Range("G28").Value = "'21"
$ws.Range("G29").Value = "'21"
$ws.Range("G30").Value = "'21"
$ws.Range("G31").Value = "'21"
$ws.Range("G32").Value = "'21"
$ws.Range("G33").Value = "'21"
$ws.Range("G34").Value = "'21"
$ws.Range("G35").Value = "'21"
$ws.Range("G36").Value = "'21"
$ws.Range("G37").Value = "'21"
$ws.Range("G38").Value = "'21"
$ws.Range("G39").Value = "'21"
$ws.Range("D40").Value = "'0.03979"
$ws.Range("G40").Value = "'21"
$ws.Range("D41").Value = "'0.006171"
$ws.Range("G41").Value = "'21"
$ws.Range("D42").Value = "'0.1070"
$ws.Range("G42").Value = "'21"
$ws.Range("G43").Value = "'21"
$ws.Range("D44").Value = "'0.009487"
$ws.Range("E44").Value = "'43LocalTradersLCTBestin24h"
$ws.Range("G44").Value = "'21"
$ws.Range("D45").Value = "'0.00005222"
$ws.Range("G45").Value = "'21"
$ws.Range("G46").Value = "'21"
$ws.Range("D47").Value = "'0.6699"
$ws.Range("G47").Value = "'21"
$ws.Range("D48").Value = "'0.002337"
$ws.Range("G48").Value = "'21"
$ws.Range("G49").Value = "'21"
$ws.Range("G50").Value = "'21"
$ws.Range("G51").Value = "'21"
